$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C7").Value = -12.8571
$ws.Range("B8").Value = 6.3755
$ws.Range("B10").Value = 5.499399999999998
$ws.Range("B12").Value = 4.582499999999998
$ws.Range("C14").Value = -12.77659999999999
$ws.Range("C15").Value = -14.00289999999998
$ws.Range("B18").Value = 7.676899999999996
$ws.Range("C18").Value = -12.16209999999999
$ws.Range("C20").Value = -11.51280000000001
$ws.Range("B25").Value = 5.907700000000001
$ws.Range("C29").Value = -11.51700000000001
$ws.Range("C30").Value = -13.0319
$ws.Range("C31").Value = -12.96
$ws.Range("C35").Value = -11.6865
$ws.Range("B37").Value = 9.297499999999998
$ws.Range("C40").Value = -13.3254
$ws.Range("C44").Value = -13.00189999999999
$ws.Range("C50").Value = -13.69109999999999
$ws.Range("C54").Value = -13.1279
$ws.Range("B55").Value = 5.796599999999995
$ws.Range("B68").Value = 6.292199999999998
$ws.Range("C68").Value = -12.12730000000001
$ws.Range("C76").Value = -12.4017
$ws.Range("B77").Value = 9.445700000000004
$ws.Range("B78").Value = 9.530600000000005
$ws.Range("B79").Value = 8.245600000000001
$ws.Range("B80").Value = 9.457699999999999
$ws.Range("B81").Value = 5.8697
$ws.Range("B82").Value = 5.445600000000001
$ws.Range("B84").Value = 6.521200000000004
$ws.Range("C87").Value = -14.08999999999999
$ws.Range("C88").Value = -12.91489999999999
$ws.Range("C92").Value = -10.87949999999999
$ws.Range("C96").Value = -12.7614
$ws.Range("C98").Value = -12.1389
$ws.Range("B101").Value = 9.209399999999993
$ws.Range("C101").Value = -13.07860000000001
$ws.Range("B102").Value = 8.309099999999999
$ws.Range("C102").Value = -12.8775
